# Apply the two changes captured in the commit:
#   1. Slide 5's table switches to a different table style (GUID change).
#   2. The theme bound to the deck's slide master (ppt/theme/theme2.xml,
#      "Red Violet"/Integral) is repainted with the plain "Office" colour
#      palette that the other theme part (ppt/theme/theme1.xml) uses.

$p = $ppt.ActivePresentation

# --- 1. Table style swap on the slide-5 table -----------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{58295048-8293-4D2C-8E5D-82B83CB8745B}")
    }
}

# --- 2. Recolour the active theme's colour scheme -------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (in that slot order).
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$firstSlide = $p.Slides.Item(1)
$themeColors = $firstSlide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColors.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
